$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.429493140873177
$ws.Range("D2").Value = 0.2206680388699453
$ws.Range("E2").Value = 0.1811656080780608
$ws.Range("F2").Value = 3.514449379536416
$ws.Range("G2").Value = 0.002605142250926588
$ws.Range("K2").Value = 4.451043915544176
$ws.Range("L2").Value = 0.1578622774669043
$ws.Range("N2").Value = 2.195141281712679
$ws.Range("C3").Value = 0.421917735229016
$ws.Range("D3").Value = 0.2219117272248994
$ws.Range("E3").Value = 0.1775845195624086
$ws.Range("F3").Value = 3.4632038256253
$ws.Range("G3").Value = 0.002612591057421148
$ws.Range("K3").Value = 4.239755081802457
$ws.Range("L3").Value = 0.1543974599136106
$ws.Range("N3").Value = 2.21424496533087
$ws.Range("C4").Value = 0.4175307019142167
$ws.Range("D4").Value = 0.2227766145306234
$ws.Range("E4").Value = 0.1755006367484526
$ws.Range("F4").Value = 3.433988170384339
$ws.Range("G4").Value = 0.00261739676370585
$ws.Range("K4").Value = 4.11305017273736
$ws.Range("L4").Value = 0.1523730292029555
$ws.Range("N4").Value = 2.226723944112841
$ws.Range("C5").Value = 0.4158089837012824
$ws.Range("D5").Value = 0.2231543894908938
$ws.Range("E5").Value = 0.174680125680549
$ws.Range("F5").Value = 3.422643444272083
$ws.Range("G5").Value = 0.002619413730798759
$ws.Range("K5").Value = 4.062170008604312
$ws.Range("L5").Value = 0.1515737432114932
$ws.Range("N5").Value = 2.231996737216257
$ws.Range("C6").Value = 0.4155270697053197
$ws.Range("D6").Value = 0.223218644461987
$ws.Range("E6").Value = 0.1745456070994997
$ws.Range("F6").Value = 3.420793401014535
$ws.Range("G6").Value = 0.002619752193263138
$ws.Range("K6").Value = 4.053766648483474
$ws.Range("L6").Value = 0.151442567709779
$ws.Range("N6").Value = 2.232883579550062
$ws.Range("C7").Value = 0.4175072154076531
$ws.Range("D7").Value = 0.2227816069668052
$ws.Range("E7").Value = 0.175489455165053
$ws.Range("F7").Value = 3.433832906581756
$ws.Range("G7").Value = 0.002617423727649444
$ws.Range("K7").Value = 4.112360946304591
$ws.Range("L7").Value = 0.152362146024096
$ws.Range("N7").Value = 2.226794296797898
$ws.Range("C8").Value = 0.4268260236963783
$ws.Range("D8").Value = 0.2210757641005898
$ws.Range("E8").Value = 0.1799068681959497
$ws.Range("F8").Value = 3.496310178777264
$ws.Range("G8").Value = 0.002607662568454418
$ws.Range("K8").Value = 4.377558198483769
$ws.Range("L8").Value = 0.1566461058381634
$ws.Range("N8").Value = 2.201572187215817
$ws.Range("C9").Value = 0.4472183274596091
$ws.Range("D9").Value = 0.2185403187370412
$ws.Range("E9").Value = 0.1894915325657109
$ws.Range("F9").Value = 3.636901408077335
$ws.Range("G9").Value = 0.002590351783445133
$ws.Range("K9").Value = 4.922034081409038
$ws.Range("L9").Value = 0.1658742857813138
$ws.Range("N9").Value = 2.158096589919793
$ws.Range("C10").Value = 0.4635237364617808
$ws.Range("D10").Value = 0.2171800643971835
$ws.Range("E10").Value = 0.1971110256650519
$ws.Range("F10").Value = 3.751544356319613
$ws.Range("G10").Value = 0.002578734397484128
$ws.Range("K10").Value = 5.33756686360806
$ws.Range("L10").Value = 0.1731739295865395
$ws.Range("N10").Value = 2.129858551138511
$ws.Range("C11").Value = 0.4712360502825277
$ws.Range("D11").Value = 0.2166723658517782
$ws.Range("E11").Value = 0.2007062656415997
$ws.Range("F11").Value = 3.806237587559906
$ws.Range("G11").Value = 0.002573685101394831
$ws.Range("K11").Value = 5.530110349364293
$ws.Range("L11").Value = 0.1766110370794678
$ws.Range("N11").Value = 2.117828648914951
$ws.Range("C12").Value = 0.4741994759384909
$ws.Range("D12").Value = 0.216496256492789
$ws.Range("E12").Value = 0.2020865282173432
$ws.Range("F12").Value = 3.827319906710954
$ws.Range("G12").Value = 0.002571806675955246
$ws.Range("K12").Value = 5.60353785998592
$ws.Range("L12").Value = 0.1779295999442496
$ws.Range("N12").Value = 2.113391600020236
$ws.Range("C13").Value = 0.4735593307915167
$ws.Range("D13").Value = 0.2165334639725103
$ws.Range("E13").Value = 0.201788422599904
$ws.Range("F13").Value = 3.822762851548333
$ws.Range("G13").Value = 0.002572209736394493
$ws.Range("K13").Value = 5.58770084484496
$ws.Range("L13").Value = 0.1776448633705314
$ws.Range("N13").Value = 2.114341913395563
$ws.Range("C14").Value = 0.4714789894084674
$ws.Range("D14").Value = 0.2166575526650192
$ws.Range("E14").Value = 0.2008194419166429
$ws.Range("F14").Value = 3.807964570532391
$ws.Range("G14").Value = 0.00257352988968271
$ws.Range("K14").Value = 5.536140886127612
$ws.Range("L14").Value = 0.176719173737041
$ws.Range("N14").Value = 2.117461230580886
$ws.Range("C15").Value = 0.4702103283578936
$ws.Range("D15").Value = 0.2167356682744384
$ws.Range("E15").Value = 0.2002283724346796
$ws.Range("F15").Value = 3.798948706097804
$ws.Range("G15").Value = 0.00257434289452152
$ws.Range("K15").Value = 5.504626351926845
$ws.Range("L15").Value = 0.1761543844041142
$ws.Range("N15").Value = 2.119387356863399
$ws.Range("C16").Value = 0.4630256980365175
$ws.Range("D16").Value = 0.2172154935641615
$ws.Range("E16").Value = 0.1968786880190407
$ws.Range("F16").Value = 3.748021652737521
$ws.Range("G16").Value = 0.002579069099684621
$ws.Range("K16").Value = 5.325055279235698
$ws.Range("L16").Value = 0.1729516719303348
$ws.Range("N16").Value = 2.130661242684525
$ws.Range("C17").Value = 0.4586940661461654
$ws.Range("D17").Value = 0.2175384232895539
$ws.Range("E17").Value = 0.1948570118061923
$ws.Range("F17").Value = 3.717434532287342
$ws.Range("G17").Value = 0.002582028628298324
$ws.Range("K17").Value = 5.2158016441565
$ws.Range("L17").Value = 0.1710169217049611
$ws.Range("N17").Value = 2.137787118199469
$ws.Range("C18").Value = 0.4562303507782701
$ws.Range("D18").Value = 0.2177346095890726
$ws.Range("E18").Value = 0.1937063321443659
$ws.Range("F18").Value = 3.700080506351412
$ws.Range("G18").Value = 0.002583753051976675
$ws.Range("K18").Value = 5.153292131209867
$ws.Range("L18").Value = 0.1699150507774903
$ws.Range("N18").Value = 2.141962445107026
$ws.Range("C19").Value = 0.4554009248301156
$ws.Range("D19").Value = 0.2178028238160437
$ws.Range("E19").Value = 0.1933188081885362
$ws.Range("F19").Value = 3.69424558691091
$ws.Range("G19").Value = 0.002584340729146574
$ws.Range("K19").Value = 5.132183913010977
$ws.Range("L19").Value = 0.1695438484996004
$ws.Range("N19").Value = 2.143389280774201
$ws.Range("C20").Value = 0.4591523025362108
$ws.Range("D20").Value = 0.217502964559543
$ws.Range("E20").Value = 0.1950709649032163
$ws.Range("F20").Value = 3.720665819739196
$ws.Range("G20").Value = 0.002581711287489167
$ws.Range("K20").Value = 5.227397625394701
$ws.Range("L20").Value = 0.1712217444314632
$ws.Range("N20").Value = 2.137020609215128
$ws.Range("C21").Value = 0.4720888664025438
$ws.Range("D21").Value = 0.216620665253707
$ws.Range("E21").Value = 0.2011035420711806
$ws.Range("F21").Value = 3.812301064580254
$ws.Range("G21").Value = 0.002573141217579454
$ws.Range("K21").Value = 5.551271225509026
$ws.Range("L21").Value = 0.1769906076283974
$ws.Range("N21").Value = 2.116541788793796
$ws.Range("C22").Value = 0.4807942030375614
$ws.Range("D22").Value = 0.2161382197259059
$ws.Range("E22").Value = 0.2051560216827468
$ws.Range("F22").Value = 3.874356251667933
$ws.Range("G22").Value = 0.002567736114425002
$ws.Range("K22").Value = 5.765950716568682
$ws.Range("L22").Value = 0.1808601349347612
$ws.Range("N22").Value = 2.103848436457369
$ws.Range("C23").Value = 0.4761248960644764
$ws.Range("D23").Value = 0.2163870342367247
$ws.Range("E23").Value = 0.2029829989024918
$ws.Range("F23").Value = 3.841036099559147
$ws.Range("G23").Value = 0.0025706030687748
$ws.Range("K23").Value = 5.651093518231164
$ws.Range("L23").Value = 0.1787857278714
$ws.Range("N23").Value = 2.110559540612883
$ws.Range("C24").Value = 0.4589450510590041
$ws.Range("D24").Value = 0.2175189626665599
$ws.Range("E24").Value = 0.1949742005683603
$ws.Range("F24").Value = 3.719204235856921
$ws.Range("G24").Value = 0.002581854685485578
$ws.Range("K24").Value = 5.222154142806858
$ws.Range("L24").Value = 0.1711291115692433
$ws.Range("N24").Value = 2.13736690329516
$ws.Range("C25").Value = 0.4414716579125582
$ws.Range("D25").Value = 0.2191386311038244
$ws.Range("E25").Value = 0.1867983287219559
$ws.Range("F25").Value = 3.596898911884693
$ws.Range("G25").Value = 0.002594840385644513
$ws.Range("K25").Value = 4.772061265310128
$ws.Range("L25").Value = 0.163287682948166
$ws.Range("N25").Value = 2.16921156424354
